# ProductBurndown.xlsx update
# - Fill in "Accomplished story points in current sprint" (D) and "Actual Hours" (E)
#   for sprints that have now completed (rows 38-55).
# - Insert three more sprint rows (56-58) with their planned hours / actuals.
# - Extend the rolling "Actual Hours" average to include the newly added sprint.
# - Move the selection to reflect where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---- Fill in actuals for already-elapsed sprints (rows 38-55) ----
$actuals = @{
    38 = @(0, 18)
    39 = @(0, 18)
    40 = @(0, 0)
    41 = @(0, 12)
    42 = @(0, 0)
    43 = @(0, 0)
    44 = @(0, 12)
    45 = @(0, 18)
    46 = @(0, 18)
    47 = @(0, 18)
    48 = @(0, 2)
    49 = @(0, 3)
    50 = @(0, 0)
    51 = @(3, 18)
    52 = @(1, 18)
}

foreach ($r in $actuals.Keys) {
    $pair = $actuals[$r]
    $ws.Cells.Item($r, 4).Value = $pair[0]
    $ws.Cells.Item($r, 5).Value = $pair[1]
}

# Rows 53-55: only the "Actual Hours" (E) column is known so far.
$ws.Cells.Item(53, 5).Value = 18
$ws.Cells.Item(54, 5).Value = 18
$ws.Cells.Item(55, 5).Value = 12

# ---- Add three more sprint rows, pushing the totals row down ----
$ws.Rows("56:58").Insert()

# Carry the formulas/number-formats of the last sprint row (55) into the new rows.
$ws.Range("A55:K55").Copy()
$ws.Range("A56:K58").PasteSpecial(-4122)

# Sprint 22 (row 56)
$ws.Cells.Item(56, 2).Value = 22
$ws.Cells.Item(56, 5).Value = 18

# Sprint 23 (row 57) - nothing accomplished/booked yet
$ws.Cells.Item(57, 2).Value = 23

# Sprint 24 (row 58)
$ws.Cells.Item(58, 2).Value = 24
$ws.Cells.Item(58, 5).Value = 2

# ---- Extend the "Actual Hours" average to cover the newly finished sprint ----
$ws.Range("E59").Formula = "=AVERAGE(E7:E56)"

# ---- Restore the user's last selection/scroll position ----
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D53").Select()

$wb.Save()
